$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.507.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.722.96'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.85%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.72'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5420'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.28%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2763'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06768'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +6.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.57'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07781'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.725'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.747.84'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.957.86'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5984'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +5.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8400'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.77'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '27.447.08'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.808'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '210.11'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +9.21%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.236'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.09%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.43'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1253'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.439'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.91'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.625'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +8.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05584'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.314'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.671'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.516'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.631'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9760'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.851'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.440'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5845'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.12%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.850'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.041.62'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8392'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.53'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.862.94'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.60'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₈109'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.21%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4405'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.9977'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05282'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.57%  '
